$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.719.07'
$ws.Range("E2").Value = '  +2.39%  '
$ws.Range("D3").Value = '2.526.20'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.97'
$ws.Range("E5").Value = '  +2.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.34'
$ws.Range("E6").Value = '  +5.65%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.14%  '
$ws.Range("D9").Value = '2.525.48'
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.141'
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.15'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.78'
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").Value = '2.985.34'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("D17").Value = '67.588.44'
$ws.Range("E17").Value = '  +2.55%  '
$ws.Range("D18").Value = '2.523.10'
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.02'
$ws.Range("E19").Value = '  +4.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.44'
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '358.90'
$ws.Range("E21").Value = '  +3.65%  '
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("E24").Value = '  +4.27%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.30'
$ws.Range("E26").Value = '  +4.69%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.82'
$ws.Range("E27").Value = '  +2.83%  '
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("E29").Value = '  -1.12%  '
$ws.Range("D30").Value = '0.0₃0988'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '555.01'
$ws.Range("E31").Value = '  +6.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.27'
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.35'
$ws.Range("E33").Value = '  +3.01%  '
$ws.Range("E34").Value = '  +2.80%  '
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  +2.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '155.65'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.74'
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.59'
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.81'
$ws.Range("E41").Value = '  +3.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.355'
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("E43").Value = '  +2.38%  '
$ws.Range("E44").Value = '  +5.71%  '
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '147.35'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("E47").Value = '  +0.81%  '
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0758'
$ws.Range("E51").Value = '  +0.24%  '
